# Weekly Fruta/Hortaliza update:
# Insert 5 new rows of "Uva" price records (week of 2022-03-08, serial 44628)
# above the previously-last rows, which shift down from 508-514 to 513-519.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 5 blank rows starting at row 508; existing rows 508:514 shift to 513:519.
$ws.Rows.Item(508).Resize(5, 1).Insert()

# Helper arrays describing the 5 new rows (columns A..T), in row order.
$newRows = @(
    @{ K = 'Flame Seedless';     L = 'Primera'; M = 380; N = 10000; O = 10000; P = 10000; Q = '$/bandeja 18 kilos'; R = 'Provincia de Los Andes'; S = 556; T = 18 },
    @{ K = 'Red Globe';          L = 'Primera'; M = 350; N = 12000; O = 12000; P = 12000; Q = '$/bandeja 18 kilos'; R = 'Provincia de Limarí';    S = 667; T = 18 },
    @{ K = 'Red Globe';          L = 'Primera'; M = 300; N = 10000; O = 10000; P = 10000; Q = '$/bandeja 18 kilos'; R = "Región de O'Higgins";    S = 556; T = 18 },
    @{ K = 'Sultanina';          L = 'Primera'; M = 350; N = 12000; O = 12000; P = 12000; Q = '$/bandeja 18 kilos'; R = "Región de O'Higgins";    S = 667; T = 18 },
    @{ K = 'Superior Seedless';  L = 'Primera'; M = 350; N = 10000; O = 10000; P = 10000; Q = '$/bandeja 18 kilos'; R = "Región de O'Higgins";    S = 556; T = 18 }
)

$startRow = 508
for ($i = 0; $i -lt $newRows.Count; $i++) {
    $r = $startRow + $i
    $row = $newRows[$i]

    $ws.Cells.Item($r, 1).Value = 9
    $ws.Cells.Item($r, 2).Value = "Vega Central Mapocho de Santiago"
    $ws.Cells.Item($r, 3).Value = "Metropolitana"
    $ws.Cells.Item($r, 4).Value2 = 44628
    $ws.Cells.Item($r, 5).Value = 13
    $ws.Cells.Item($r, 6).Value = "Fruta"
    $ws.Cells.Item($r, 7).Value = 100109
    $ws.Cells.Item($r, 8).Value = "Uva"
    $ws.Cells.Item($r, 9).Value = 100109001
    $ws.Cells.Item($r, 10).Value = "Uva"
    $ws.Cells.Item($r, 11).Value = $row.K
    $ws.Cells.Item($r, 12).Value = $row.L
    $ws.Cells.Item($r, 13).Value = $row.M
    $ws.Cells.Item($r, 14).Value = $row.N
    $ws.Cells.Item($r, 15).Value = $row.O
    $ws.Cells.Item($r, 16).Value = $row.P
    $ws.Cells.Item($r, 17).Value = $row.Q
    $ws.Cells.Item($r, 18).Value = $row.R
    $ws.Cells.Item($r, 19).Value = $row.S
    $ws.Cells.Item($r, 20).Value = $row.T
}
